# "update number of cities"
#
# The author refreshed the population figures for Buenos Aires and Toronto,
# then re-sorted the city table by Population (column C) descending, and
# turned on AutoFilter for the table range A1:D13 (which is also what makes
# Excel register the hidden "_FilterDatabase" defined name for the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated population counts (row positions are still in their original,
# pre-sort order at this point).
$ws.Range("C4").Value = 2891000    # Buenos Aires, Argentina
$ws.Range("C12").Value = 2800000   # Toronto, Canada

# Re-sort the table A1:D13 by Population (column C), largest to smallest.
$table = $ws.Range("A1:D13")
$table.Sort($ws.Range("C1"), 2)

# Turn on AutoFilter for the (now sorted) table.
$ws.Range("A1:D13").AutoFilter()

# Excel records the filtered range as a hidden, sheet-scoped defined name.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "Sheet1!`$A`$1:`$D`$13")
$filterName.Visible = $false
